$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set progress value for "Trần Quốc Thịnh" task 1 to 100% in C3
$ws.Range("C3").Value = 1
$ws.Range("C3").NumberFormat = "0%"

# Move the active selection (cosmetic, matches the authored selection state)
$ws.Range("F5").Select()
